$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 header: "pt_min" -> "eta" (column F no longer holds a constant pt_min,
# it now holds per-row eta values)
$ws.Range("F1").Value = "eta"

# Column F (rows 2-15): replace the constant 25 (pt_min) with per-row eta
# values for each eta bin.
$ws.Range("F2").Value = 0.1
$ws.Range("F3").Value = 0.3
$ws.Range("F4").Value = 0.5
$ws.Range("F5").Value = 0.7
$ws.Range("F6").Value = 0.9
$ws.Range("F7").Value = 1.1000000000000001
$ws.Range("F8").Value = 1.3
$ws.Range("F9").Value = 1.5
$ws.Range("F10").Value = 1.7
$ws.Range("F11").Value = 1.9
$ws.Range("F12").Value = 2.1
$ws.Range("F13").Value = 2.29
$ws.Range("F14").Value = 2.52
$ws.Range("F15").Value = 2.81

# Move the active selection from H19 to G19.
$ws.Range("G19").Select()
